$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as Text so numeric-looking strings
# (e.g. "569.66") are not auto-converted to numbers by Excel.
$colD = $ws.Range("D2:D51")
$colD.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "64.055.40"
$ws.Range("E2").Value = "  +0.03%  "

# Row 3
$ws.Range("D3").Value = "2.739.02"
$ws.Range("E3").Value = "  -0.56%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").Value = "569.66"
$ws.Range("E5").Value = "  -1.30%  "

# Row 6
$ws.Range("D6").Value = "158.95"
$ws.Range("E6").Value = "  +0.43%  "

# Row 7
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$ws.Range("D8").Value = "0.598"
$ws.Range("E8").Value = "  -1.62%  "

# Row 9
$ws.Range("E9").Value = "  -1.62%  "

# Row 10
$ws.Range("E10").Value = "  +4.38%  "

# Row 11
$ws.Range("E11").Value = "  -2.48%  "

# Row 12
$ws.Range("D12").Value = "0.383"
$ws.Range("E12").Value = "  -1.06%  "

# Row 13
$ws.Range("D13").Value = "3.221.07"
$ws.Range("E13").Value = "  -0.70%  "

# Row 14
$ws.Range("D14").Value = "26.69"
$ws.Range("E14").Value = "  -0.70%  "

# Row 15
$ws.Range("D15").Value = "63.631.42"
$ws.Range("E15").Value = "  -0.46%  "

# Row 16
$ws.Range("E16").Value = "  -1.89%  "

# Row 17
$ws.Range("D17").Value = "2.743.49"
$ws.Range("E17").Value = "  -0.56%  "

# Row 18
$ws.Range("D18").Value = "12.08"
$ws.Range("E18").Value = "  -0.31%  "

# Row 19
$ws.Range("E19").Value = "  -1.89%  "

# Row 20
$ws.Range("D20").Value = "354.46"
$ws.Range("E20").Value = "  -1.69%  "

# Row 21
$ws.Range("D21").Value = "6.60"
$ws.Range("E21").Value = "  -3.20%  "

# Row 22
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.02%  "

# Row 23
$ws.Range("D23").Value = "0.523"
$ws.Range("E23").Value = "  -5.64%  "

# Row 24
$ws.Range("D24").Value = "64.35"
$ws.Range("E24").Value = "  -2.89%  "

# Row 25
$ws.Range("E25").Value = "  -0.17%  "

# Row 26
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  +0.15%  "

# Row 27
$ws.Range("E27").Value = "  -1.06%  "

# Row 28
$ws.Range("E28").Value = "  -2.62%  "

# Row 29
$ws.Range("E29").Value = "  -0.44%  "

# Row 30
$ws.Range("D30").Value = "7.23"
$ws.Range("E30").Value = "  +2.65%  "

# Row 31
$ws.Range("E31").Value = "  +7.56%  "

# Row 32
$ws.Range("D32").Value = "164.27"
$ws.Range("E32").Value = "  -3.26%  "

# Row 33
$ws.Range("E33").Value = "  -0.91%  "

# Row 34
$ws.Range("D34").Value = "20.02"
$ws.Range("E34").Value = "  -1.83%  "

# Row 35
$ws.Range("E35").Value = "  +1.33%  "

# Row 37
$ws.Range("E37").Value = "  +0.28%  "

# Row 38
$ws.Range("D38").Value = "0.991"
$ws.Range("E38").Value = "  -1.06%  "

# Row 39
$ws.Range("D39").Value = "350.10"
$ws.Range("E39").Value = "  +5.65%  "

# Row 40
$ws.Range("E40").Value = "  +1.62%  "

# Row 41
$ws.Range("E41").Value = "  -1.25%  "

# Row 42
$ws.Range("D42").Value = "38.57"
$ws.Range("E42").Value = "  -1.62%  "

# Row 43
$ws.Range("D43").Value = "22.05"
$ws.Range("E43").Value = "  +1.15%  "

# Row 44
$ws.Range("D44").Value = "21.15"
$ws.Range("E44").Value = "  -3.29%  "

# Row 45
$ws.Range("D45").Value = "0.0584"
$ws.Range("E45").Value = "  -1.93%  "

# Row 46
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").Value = "0.626"
$ws.Range("E46").Value = "  -1.67%  "

# Row 47
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "134.70"
$ws.Range("E47").Value = "  -1.13%  "

# Row 48
$ws.Range("E48").Value = "  -1.24%  "

# Row 49
$ws.Range("E49").Value = "  -3.54%  "

# Row 50
$ws.Range("E50").Value = "  -0.11%  "

# Row 51
$ws.Range("D51").Value = "11.05"
$ws.Range("E51").Value = "  +0.02%  "

# Restore column D to its original (default/Normal) style so we don't
# leave a stray number-format override behind.
$colD.Style = "Normal"